$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 361 (old rows 361..464 shift down to 363..466)
$ws.Rows("361:362").Insert()

# Populate new row 361
$ws.Cells.Item(361, 1).Value = 10
$ws.Cells.Item(361, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(361, 3).Value = "La Araucanía"
$ws.Cells.Item(361, 4).Value = 44551
$ws.Cells.Item(361, 5).Value = 9
$ws.Cells.Item(361, 6).Value = 100112003
$ws.Cells.Item(361, 7).Value = "Ajo"
$ws.Cells.Item(361, 8).Value = "Chino"
$ws.Cells.Item(361, 9).Value = "Primera"
$ws.Cells.Item(361, 10).Value = 235
$ws.Cells.Item(361, 11).Value = 20000
$ws.Cells.Item(361, 12).Value = 21000
$ws.Cells.Item(361, 13).Value = 20532
$ws.Cells.Item(361, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(361, 15).Value = "China"
$ws.Cells.Item(361, 16).Value = 2053
$ws.Cells.Item(361, 17).Value = 10
$ws.Cells.Item(361, 18).Value = "Hortaliza"

# Populate new row 362
$ws.Cells.Item(362, 1).Value = 10
$ws.Cells.Item(362, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(362, 3).Value = "La Araucanía"
$ws.Cells.Item(362, 4).Value = 44551
$ws.Cells.Item(362, 5).Value = 9
$ws.Cells.Item(362, 6).Value = 100112003
$ws.Cells.Item(362, 7).Value = "Ajo"
$ws.Cells.Item(362, 8).Value = "Chino"
$ws.Cells.Item(362, 9).Value = "Primera"
$ws.Cells.Item(362, 10).Value = 110
$ws.Cells.Item(362, 11).Value = 23000
$ws.Cells.Item(362, 12).Value = 23000
$ws.Cells.Item(362, 13).Value = 23000
$ws.Cells.Item(362, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(362, 15).Value = "China"
$ws.Cells.Item(362, 16).Value = 2300
$ws.Cells.Item(362, 17).Value = 10
$ws.Cells.Item(362, 18).Value = "Hortaliza"

# Ensure D column on new rows keeps the date number format (style index 2), matching the
# style already used throughout column D.
$ws.Cells.Item(361, 4).NumberFormat = $ws.Cells.Item(360, 4).NumberFormat
$ws.Cells.Item(362, 4).NumberFormat = $ws.Cells.Item(360, 4).NumberFormat
